$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1137-1138, pushing all existing data (old rows
# 1137..1227) down by two rows (they become 1139..1229). This mirrors the
# canonical diff, which shows a new week of price data being inserted
# right after the most-recent existing entries (rows 1135/1136) and the
# sheet's dimension growing from A1:T1227 to A1:T1229.
$ws.Rows("1137:1138").Insert()

# New row 1137: "Pintón" quality entry for the newly-added date.
$ws.Cells.Item(1137, 1).Value = 5
$ws.Cells.Item(1137, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1137, 3).Value = "Maule"
$ws.Cells.Item(1137, 4).Value = 45223
$ws.Cells.Item(1137, 5).Value = 7
$ws.Cells.Item(1137, 6).Value = "Fruta"
$ws.Cells.Item(1137, 7).Value = 100108
$ws.Cells.Item(1137, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1137, 9).Value = 100108006
$ws.Cells.Item(1137, 10).Value = "Plátano"
$ws.Cells.Item(1137, 11).Value = "Sin especificar"
$ws.Cells.Item(1137, 12).Value = "Pintón"
$ws.Cells.Item(1137, 13).Value = 1000
$ws.Cells.Item(1137, 14).Value = 22000
$ws.Cells.Item(1137, 15).Value = 22000
$ws.Cells.Item(1137, 16).Value = 22000
$ws.Cells.Item(1137, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1137, 18).Value = "Ecuador"
$ws.Cells.Item(1137, 19).Value = 1100
$ws.Cells.Item(1137, 20).Value = 20

# New row 1138: "Primera Pintón" quality entry for the same new date.
$ws.Cells.Item(1138, 1).Value = 5
$ws.Cells.Item(1138, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1138, 3).Value = "Maule"
$ws.Cells.Item(1138, 4).Value = 45223
$ws.Cells.Item(1138, 5).Value = 7
$ws.Cells.Item(1138, 6).Value = "Fruta"
$ws.Cells.Item(1138, 7).Value = 100108
$ws.Cells.Item(1138, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1138, 9).Value = 100108006
$ws.Cells.Item(1138, 10).Value = "Plátano"
$ws.Cells.Item(1138, 11).Value = "Sin especificar"
$ws.Cells.Item(1138, 12).Value = "Primera Pintón"
$ws.Cells.Item(1138, 13).Value = 600
$ws.Cells.Item(1138, 14).Value = 23000
$ws.Cells.Item(1138, 15).Value = 24000
$ws.Cells.Item(1138, 16).Value = 23333
$ws.Cells.Item(1138, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1138, 18).Value = "Ecuador"
$ws.Cells.Item(1138, 19).Value = 1167
$ws.Cells.Item(1138, 20).Value = 20
